$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append three new rows (32-34) of data, continuing the A/B/C pattern
# found in the preceding block of rows.
$data = @(
    @(1, 9.138900628687384, 1.34558127125136),
    @(2, 8.71032033017539, 1.480102726150887),
    @(3, 9.898999999999999, 1.081122146288446)
)

$startRow = 32
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
